$wb = $excel.ActiveWorkbook

# Event sheet: add event__series, event__city, event__country, event__region
$wsEvent = $wb.Worksheets.Item("Event")
$wsEvent.Range("F1").Value = "event__series"
$wsEvent.Range("G1").Value = "event__city"
$wsEvent.Range("H1").Value = "event__country"
$wsEvent.Range("I1").Value = "event__region"

# Institution sheet: add institution__scholar
$wsInstitution = $wb.Worksheets.Item("Institution")
$wsInstitution.Range("E1").Value = "institution__scholar"
